$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.219905376434326
$ws.Range("B1").Value = 3.103070259094238
$ws.Range("C1").Value = 2.740068912506104
$ws.Range("D1").Value = 2.46838116645813
$ws.Range("E1").Value = 1.727951884269714
